# Update countries & provincias Spain
# Applies the 9-Sept-2020 19:36 COVID data refresh to the "Pais" sheet:
#  - Refreshes the "Datos actualizados..." timestamp in A1
#  - Updates Casos totales/Nuevos casos/Casos activos/Recuperados/Casos
#    criticos/Muertes hoy/Muertes for the countries whose stats moved
#  - Republica Dominicana overtakes Egipto in the ranking (rows 34/35 swap
#    country labels), and Libano overtakes Corea del Sur (rows 78/79 swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 9 de Septiembre de 2020 a las 19:36
$ws.Range("A1").Value2 = 'Datos actualizados a 9 de Septiembre de 2020 a las 19:36'

# Row 4: Estados Unidos
$ws.Range("B4").Value2 = 6526716
$ws.Range("C4").Value2 = 12485
$ws.Range("D4").Value2 = 3807163
$ws.Range("E4").Value2 = 2525090
$ws.Range("G4").Value2 = 433
$ws.Range("H4").Value2 = 194463

# Row 5: India
$ws.Range("B5").Value2 = 4433708
$ws.Range("C5").Value2 = 66272
$ws.Range("D5").Value2 = 3449450
$ws.Range("E5").Value2 = 909645
$ws.Range("G5").Value2 = 690
$ws.Range("H5").Value2 = 74613

# Row 12: España
$ws.Range("B12").Value2 = 543379
$ws.Range("C12").Value2 = 8866
$ws.Range("G12").Value2 = 34
$ws.Range("H12").Value2 = 29628

# Row 24: Alemania
$ws.Range("B24").Value2 = 255681
$ws.Range("C24").Value2 = 725
$ws.Range("E24").Value2 = 15672

# Row 28: Israel
$ws.Range("B28").Value2 = 140658
$ws.Range("C28").Value2 = 3093
$ws.Range("D28").Value2 = 108339
$ws.Range("E28").Value2 = 31266
$ws.Range("G28").Value2 = 13
$ws.Range("H28").Value2 = 1053

# Row 29: Canada
$ws.Range("B29").Value2 = 134077
$ws.Range("C29").Value2 = 329
$ws.Range("D29").Value2 = 117939
$ws.Range("E29").Value2 = 6984
$ws.Range("G29").Value2 = 1
$ws.Range("H29").Value2 = 9154

# Row 32: Ecuador
$ws.Range("B32").Value2 = 112166
$ws.Range("C32").Value2 = 1409
$ws.Range("E32").Value2 = 10223
$ws.Range("G32").Value2 = 74
$ws.Range("H32").Value2 = 10701

# Row 34: Republica Dominicana
$ws.Range("A34").Value2 = 'Republica Dominicana'
$ws.Range("B34").Value2 = 100937
$ws.Range("C34").Value2 = 806
$ws.Range("D34").Value2 = 74305
$ws.Range("E34").Value2 = 24718
$ws.Range("G34").Value2 = 25
$ws.Range("H34").Value2 = 1914

# Row 35: Egipto
$ws.Range("A35").Value2 = 'Egipto'
$ws.Range("B35").Value2 = 100228
$ws.Range("D35").Value2 = 79886
$ws.Range("E35").Value2 = 14782
$ws.Range("H35").Value2 = 5560

# Row 59: Argelia
$ws.Range("B59").Value2 = 47216
$ws.Range("C59").Value2 = 278
$ws.Range("D59").Value2 = 33379
$ws.Range("E59").Value2 = 12256
$ws.Range("G59").Value2 = 10
$ws.Range("H59").Value2 = 1581

# Row 71: Chequia
$ws.Range("D71").Value2 = 20365
$ws.Range("E71").Value2 = 9767

# Row 72: Irlanda
$ws.Range("B72").Value2 = 30164
$ws.Range("C72").Value2 = 84
$ws.Range("E72").Value2 = 5019
$ws.Range("G72").Value2 = 3
$ws.Range("H72").Value2 = 1781

# Row 78: Libano
$ws.Range("A78").Value2 = 'Libano'
$ws.Range("B78").Value2 = 21877
$ws.Range("C78").Value2 = 553
$ws.Range("D78").Value2 = 7024
$ws.Range("E78").Value2 = 14641
$ws.Range("G78").Value2 = 5
$ws.Range("H78").Value2 = 212

# Row 79: Corea del Sur
$ws.Range("A79").Value2 = 'Corea del Sur'
$ws.Range("B79").Value2 = 21588
$ws.Range("C79").Value2 = 156
$ws.Range("D79").Value2 = 17023
$ws.Range("E79").Value2 = 4221
$ws.Range("G79").Value2 = 3
$ws.Range("H79").Value2 = 344

# Row 104: Haiti
$ws.Range("B104").Value2 = 8384
$ws.Range("C104").Value2 = 8
$ws.Range("E104").Value2 = 2179

# Row 109: Malaui
$ws.Range("B109").Value2 = 5653
$ws.Range("C109").Value2 = 23
$ws.Range("E109").Value2 = 1847

# Row 111: Republica de Yibuti
$ws.Range("B111").Value2 = 5391
$ws.Range("C111").Value2 = 3
$ws.Range("E111").Value2 = 3

# Row 133: Sri Lanka
$ws.Range("B133").Value2 = 3147
$ws.Range("C133").Value2 = 7
$ws.Range("E133").Value2 = 189

# Row 148: Sierra Leona
$ws.Range("B148").Value2 = 2067
$ws.Range("C148").Value2 = 3
$ws.Range("D148").Value2 = 1622
$ws.Range("E148").Value2 = 373

# Row 190: Bermudas
$ws.Range("B190").Value2 = 177
$ws.Range("C190").Value2 = 2
$ws.Range("D190").Value2 = 159
$ws.Range("E190").Value2 = 9

